$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 212.4076363333333
$ws.Range("H2").Value = 637.222909
$ws.Range("I2").Value = 0.5080632835800084
$ws.Range("J2").Value = 0.5080632835800084
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 153.5290173333333
$ws.Range("N2").Value = 460.587052
$ws.Range("O2").Value = 0.3172206968818489
$ws.Range("P2").Value = 0.317220696881849
$ws.Range("Q2").Value = 32610.73568035269
$ws.Range("R2").Value = 293496.6211231742
$ws.Range("S2").Value = 0.1611681888773307
$ws.Range("T2").Value = 0.1611681888773307

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 212.4076363333333
$ws.Range("H3").Value = 637.222909
$ws.Range("I3").Value = 0.5080632835800084
$ws.Range("J3").Value = 0.5080632835800084
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 168.7997026666667
$ws.Range("N3").Value = 506.3991080000001
$ws.Range("O3").Value = 0.3487728915577651
$ws.Range("P3").Value = 0.3487728915577651
$ws.Range("Q3").Value = 35854.34585719613
$ws.Range("R3").Value = 322689.1127147652
$ws.Range("S3").Value = 0.1771987005085323
$ws.Range("T3").Value = 0.1771987005085323

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 212.4076363333333
$ws.Range("H4").Value = 637.222909
$ws.Range("I4").Value = 0.5080632835800084
$ws.Range("J4").Value = 0.5080632835800084
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 68.09032333333333
$ws.Range("N4").Value = 204.27097
$ws.Range("O4").Value = 0.1406878008722904
$ws.Range("P4").Value = 0.1406878008722904
$ws.Range("Q4").Value = 14462.90463640575
$ws.Range("R4").Value = 130166.1417276517
$ws.Range("S4").Value = 0.07147830607082623
$ws.Range("T4").Value = 0.07147830607082624

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 212.4076363333333
$ws.Range("H5").Value = 637.222909
$ws.Range("I5").Value = 0.5080632835800084
$ws.Range("J5").Value = 0.5080632835800084
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 93.562673
$ws.Range("N5").Value = 280.688019
$ws.Range("O5").Value = 0.1933186106880956
$ws.Range("P5").Value = 0.1933186106880956
$ws.Range("Q5").Value = 19873.42622095859
$ws.Range("R5").Value = 178860.8359886273
$ws.Range("S5").Value = 0.09821808812331918
$ws.Range("T5").Value = 0.09821808812331918

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 161.9384456666667
$ws.Range("H6").Value = 485.815337
$ws.Range("I6").Value = 0.3873447295187379
$ws.Range("J6").Value = 0.3873447295187379
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 153.5290173333333
$ws.Range("N6").Value = 460.587052
$ws.Range("O6").Value = 0.3172206968818489
$ws.Range("P6").Value = 0.317220696881849
$ws.Range("Q6").Value = 24862.25043169072
$ws.Range("R6").Value = 223760.2538852165
$ws.Range("S6").Value = 0.1228737650314453
$ws.Range("T6").Value = 0.1228737650314453

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 161.9384456666667
$ws.Range("H7").Value = 485.815337
$ws.Range("I7").Value = 0.3873447295187379
$ws.Range("J7").Value = 0.3873447295187379
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 168.7997026666667
$ws.Range("N7").Value = 506.3991080000001
$ws.Range("O7").Value = 0.3487728915577651
$ws.Range("P7").Value = 0.3487728915577651
$ws.Range("Q7").Value = 27335.16147883549
$ws.Range("R7").Value = 246016.4533095194
$ws.Range("S7").Value = 0.1350953413439106
$ws.Range("T7").Value = 0.1350953413439106

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 161.9384456666667
$ws.Range("H8").Value = 485.815337
$ws.Range("I8").Value = 0.3873447295187379
$ws.Range("J8").Value = 0.3873447295187379
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 68.09032333333333
$ws.Range("N8").Value = 204.27097
$ws.Range("O8").Value = 0.1406878008722904
$ws.Range("P8").Value = 0.1406878008722904
$ws.Range("Q8").Value = 11026.44112554076
$ws.Range("R8").Value = 99237.9701298669
$ws.Range("S8").Value = 0.05449467817546338
$ws.Range("T8").Value = 0.05449467817546339

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 161.9384456666667
$ws.Range("H9").Value = 485.815337
$ws.Range("I9").Value = 0.3873447295187379
$ws.Range("J9").Value = 0.3873447295187379
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 93.562673
$ws.Range("N9").Value = 280.688019
$ws.Range("O9").Value = 0.1933186106880956
$ws.Range("P9").Value = 0.1933186106880956
$ws.Range("Q9").Value = 15151.3938380386
$ws.Range("R9").Value = 136362.5445423474
$ws.Range("S9").Value = 0.07488094496791861
$ws.Range("T9").Value = 0.07488094496791861

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.338549
$ws.Range("H10").Value = 1.015647
$ws.Range("I10").Value = 0.0008097840527861261
$ws.Range("J10").Value = 0.0008097840527861261
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 153.5290173333333
$ws.Range("N10").Value = 460.587052
$ws.Range("O10").Value = 0.3172206968818489
$ws.Range("P10").Value = 0.317220696881849
$ws.Range("Q10").Value = 51.97709528918266
$ws.Range("R10").Value = 467.793857602644
$ws.Range("S10").Value = 0.0002568802615486228
$ws.Range("T10").Value = 0.0002568802615486229

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 0.338549
$ws.Range("H11").Value = 1.015647
$ws.Range("I11").Value = 0.0008097840527861261
$ws.Range("J11").Value = 0.0008097840527861261
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 168.7997026666667
$ws.Range("N11").Value = 506.3991080000001
$ws.Range("O11").Value = 0.3487728915577651
$ws.Range("P11").Value = 0.3487728915577651
$ws.Range("Q11").Value = 57.14697053809734
$ws.Range("R11").Value = 514.322734842876
$ws.Range("S11").Value = 0.000282430725627583
$ws.Range("T11").Value = 0.000282430725627583

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 0.338549
$ws.Range("H12").Value = 1.015647
$ws.Range("I12").Value = 0.0008097840527861261
$ws.Range("J12").Value = 0.0008097840527861261
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 68.09032333333333
$ws.Range("N12").Value = 204.27097
$ws.Range("O12").Value = 0.1406878008722904
$ws.Range("P12").Value = 0.1406878008722904
$ws.Range("Q12").Value = 23.05191087417667
$ws.Range("R12").Value = 207.46719786759
$ws.Range("S12").Value = 0.0001139267375679308
$ws.Range("T12").Value = 0.0001139267375679308

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 0.338549
$ws.Range("H13").Value = 1.015647
$ws.Range("I13").Value = 0.0008097840527861261
$ws.Range("J13").Value = 0.0008097840527861261
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 93.562673
$ws.Range("N13").Value = 280.688019
$ws.Range("O13").Value = 0.1933186106880956
$ws.Range("P13").Value = 0.1933186106880956
$ws.Range("Q13").Value = 31.675549381477
$ws.Range("R13").Value = 285.079944433293
$ws.Range("S13").Value = 0.0001565463280419894
$ws.Range("T13").Value = 0.0001565463280419894

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 43.38855633333333
$ws.Range("H14").Value = 130.165669
$ws.Range("I14").Value = 0.1037822028484675
$ws.Range("J14").Value = 0.1037822028484675
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 153.5290173333333
$ws.Range("N14").Value = 460.587052
$ws.Range("O14").Value = 0.3172206968818489
$ws.Range("P14").Value = 0.317220696881849
$ws.Range("Q14").Value = 6661.402417368642
$ws.Range("R14").Value = 59952.62175631779
$ws.Range("S14").Value = 0.03292186271152425
$ws.Range("T14").Value = 0.03292186271152425

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 43.38855633333333
$ws.Range("H15").Value = 130.165669
$ws.Range("I15").Value = 0.1037822028484675
$ws.Range("J15").Value = 0.1037822028484675
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 168.7997026666667
$ws.Range("N15").Value = 506.3991080000001
$ws.Range("O15").Value = 0.3487728915577651
$ws.Range("P15").Value = 0.3487728915577651
$ws.Range("Q15").Value = 7323.975408202585
$ws.Range("R15").Value = 65915.77867382327
$ws.Range("S15").Value = 0.03619641897969452
$ws.Range("T15").Value = 0.03619641897969452

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 43.38855633333333
$ws.Range("H16").Value = 130.165669
$ws.Range("I16").Value = 0.1037822028484675
$ws.Range("J16").Value = 0.1037822028484675
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 68.09032333333333
$ws.Range("N16").Value = 204.27097
$ws.Range("O16").Value = 0.1406878008722904
$ws.Range("P16").Value = 0.1406878008722904
$ws.Range("Q16").Value = 2954.340829703214
$ws.Range("R16").Value = 26589.06746732893
$ws.Range("S16").Value = 0.01460088988843284
$ws.Range("T16").Value = 0.01460088988843284

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 43.38855633333333
$ws.Range("H17").Value = 130.165669
$ws.Range("I17").Value = 0.1037822028484675
$ws.Range("J17").Value = 0.1037822028484675
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 93.562673
$ws.Range("N17").Value = 280.688019
$ws.Range("O17").Value = 0.1933186106880956
$ws.Range("P17").Value = 0.1933186106880956
$ws.Range("Q17").Value = 4059.549308157746
$ws.Range("R17").Value = 36535.94377341971
$ws.Range("S17").Value = 0.02006303126881585
$ws.Range("T17").Value = 0.02006303126881585

